$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Actualiza base de datos EC ---
# Row 17: CHRISTIAN ALFONSO MARTINEZ CASTILLO, periodo 2104 -> 2002 (ajuste de periodo/valor)
$ws.Range("E17").Value = "2002"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 877803

# Rows 18-24: ahora pertenecen a HERNAN GREGORIO BARRIOS AGUILERA (periodos 2002-2008)
$ws.Range("C18").Value = "1235048851"
$ws.Range("D18").Value = "HERNAN GREGORIO BARRIOS AGUILERA"
$ws.Range("E18").Value = "2002"
$ws.Range("F18").Value = 33200
$ws.Range("G18").Value = 830000

$ws.Range("C19").Value = "1235048851"
$ws.Range("D19").Value = "HERNAN GREGORIO BARRIOS AGUILERA"
$ws.Range("E19").Value = "2003"
$ws.Range("F19").Value = 33200
$ws.Range("G19").Value = 830000

$ws.Range("C20").Value = "1235048851"
$ws.Range("D20").Value = "HERNAN GREGORIO BARRIOS AGUILERA"
$ws.Range("E20").Value = "2004"
$ws.Range("F20").Value = 33200
$ws.Range("G20").Value = 830000

$ws.Range("C21").Value = "1235048851"
$ws.Range("D21").Value = "HERNAN GREGORIO BARRIOS AGUILERA"
$ws.Range("E21").Value = "2005"
$ws.Range("F21").Value = 33200
$ws.Range("G21").Value = 830000

$ws.Range("C22").Value = "1235048851"
$ws.Range("D22").Value = "HERNAN GREGORIO BARRIOS AGUILERA"
$ws.Range("E22").Value = "2006"
$ws.Range("F22").Value = 33200
$ws.Range("G22").Value = 830000

$ws.Range("C23").Value = "1235048851"
$ws.Range("D23").Value = "HERNAN GREGORIO BARRIOS AGUILERA"
$ws.Range("E23").Value = "2007"
$ws.Range("F23").Value = 33200
$ws.Range("G23").Value = 830000

$ws.Range("C24").Value = "1235048851"
$ws.Range("D24").Value = "HERNAN GREGORIO BARRIOS AGUILERA"
$ws.Range("E24").Value = "2008"
$ws.Range("F24").Value = 33200
$ws.Range("G24").Value = 830000

# Row 25: CHRISTIAN ALFONSO MARTINEZ CASTILLO periodo 2009
$ws.Range("E25").Value = "2009"
$ws.Range("F25").Value = 35112
$ws.Range("G25").Value = 877803

# Row 26: HERNAN GREGORIO BARRIOS AGUILERA periodo 2009
$ws.Range("E26").Value = "2009"
$ws.Range("F26").Value = 33200

# Rows 27-40: intercalado CHRISTIAN / HERNAN por periodo 2010-2104
$ws.Range("C27").Value = "73143267"
$ws.Range("D27").Value = "CHRISTIAN ALFONSO MARTINEZ CASTILLO"
$ws.Range("E27").Value = "2010"
$ws.Range("F27").Value = 35112
$ws.Range("G27").Value = 877803

$ws.Range("E28").Value = "2010"

$ws.Range("C29").Value = "73143267"
$ws.Range("D29").Value = "CHRISTIAN ALFONSO MARTINEZ CASTILLO"
$ws.Range("E29").Value = "2011"
$ws.Range("F29").Value = 35112
$ws.Range("G29").Value = 877803

$ws.Range("E30").Value = "2011"

$ws.Range("C31").Value = "73143267"
$ws.Range("D31").Value = "CHRISTIAN ALFONSO MARTINEZ CASTILLO"
$ws.Range("E31").Value = "2012"
$ws.Range("F31").Value = 35112
$ws.Range("G31").Value = 877803

$ws.Range("E32").Value = "2012"

$ws.Range("C33").Value = "73143267"
$ws.Range("D33").Value = "CHRISTIAN ALFONSO MARTINEZ CASTILLO"
$ws.Range("E33").Value = "2101"
$ws.Range("F33").Value = 35112
$ws.Range("G33").Value = 877803

$ws.Range("E34").Value = "2101"

$ws.Range("C35").Value = "73143267"
$ws.Range("D35").Value = "CHRISTIAN ALFONSO MARTINEZ CASTILLO"
$ws.Range("E35").Value = "2102"
$ws.Range("F35").Value = 35112
$ws.Range("G35").Value = 877803

$ws.Range("E36").Value = "2102"

$ws.Range("C37").Value = "73143267"
$ws.Range("D37").Value = "CHRISTIAN ALFONSO MARTINEZ CASTILLO"
$ws.Range("E37").Value = "2103"
$ws.Range("F37").Value = 35112
$ws.Range("G37").Value = 877803

$ws.Range("E38").Value = "2103"

$ws.Range("C39").Value = "73143267"
$ws.Range("D39").Value = "CHRISTIAN ALFONSO MARTINEZ CASTILLO"
$ws.Range("E39").Value = "2104"
$ws.Range("F39").Value = 30430
$ws.Range("G39").Value = 877803

$ws.Range("E40").Value = "2104"
$ws.Range("F40").Value = 28774
